# --- Edit: "modified vendors' label and switch" ---------------------------
# The author filled in every previously-blank cell across the used range
# (rows 2-27, including rows 3, 6, 8, 10, 13 and 17 which were entirely
# blank) with an empty text entry, and appended a brand-new order on row 29
# (which extends the sheet's used range from A1:AM28 to A1:AM29).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> list of column letters that are currently empty in
# that row and need an empty-string text entry.
$blankCellsByRow = @{
    2 = @("L", "Q", "R", "S", "T", "AH", "AI", "AJ", "AK", "AL", "AM")
    3 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    4 = @("L", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    5 = @("L", "V", "W", "Y", "Z", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    6 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    7 = @("L", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    8 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    9 = @("L", "O", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    10 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    11 = @("L", "Q", "R", "S", "T", "V", "W", "Y", "Z", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    12 = @("L", "Q", "R", "S", "T", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    13 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    14 = @("L", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    15 = @("L", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    16 = @("L", "Q", "R", "S", "T", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    17 = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T", "U", "V", "W", "X", "Y", "Z", "AA", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    18 = @("L", "Q", "R", "S", "T", "AH", "AI", "AJ", "AK", "AL", "AM")
    19 = @("L", "AH", "AI", "AJ", "AK", "AL", "AM")
    20 = @("L", "Q", "R", "S", "T")
    21 = @("L", "Q", "R", "S", "T", "Y", "Z", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    22 = @("AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    23 = @("Q", "R", "S", "T", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
    24 = @("Q", "R", "S", "T", "AK", "AL", "AM")
    25 = @("Q", "R", "S", "T", "AK", "AL", "AM")
    26 = @("Q", "R", "S", "T", "AH", "AI", "AJ", "AK", "AL", "AM")
    27 = @("Q", "R", "S", "T", "AK", "AL", "AM")
}

foreach ($r in $blankCellsByRow.Keys) {
    foreach ($c in $blankCellsByRow[$r]) {
        $cell = $ws.Range("$c$r")
        # A lone apostrophe enters an explicit (empty) text value without
        # Excel re-interpreting it as a number/date, matching the empty
        # inline-string cells from the diff; Style reset drops the quote-
        # prefix formatting the apostrophe trick leaves behind.
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
}

# --- New row 29 --------------------------------------------------------
$row = 29

# DATE column: force text so "03/14/2024" is stored as a literal string
# (matching every other DATE cell in this sheet) instead of being parsed
# into a date serial number.
$dateCell = $ws.Range("A$row")
$dateCell.Value = "'03/14/2024"
$dateCell.Style = "Normal"

$ws.Range("B$row").Value = "SO240314001"
$ws.Range("C$row").Value = "ab"
$ws.Range("D$row").Value = "abs@abc.com"
$ws.Range("E$row").Value = "(789)456-1233"
$ws.Range("F$row").Value = "YES"
$ws.Range("G$row").Value = "YES"
$ws.Range("H$row").Value = "artist"
$ws.Range("I$row").Value = "title"
$ws.Range("J$row").Value = 10
$ws.Range("K$row").Value = 85
$ws.Range("L$row").Value = " "
$ws.Range("M$row").Value = "Ebay"
$ws.Range("N$row").Value = "LP"
$ws.Range("O$row").Value = "abake"
$ws.Range("P$row").Value = "PICKUP"
$ws.Range("Q$row").Value = "N/A"
$ws.Range("R$row").Value = "N/A"
$ws.Range("S$row").Value = "N/A"
$ws.Range("T$row").Value = "N/A"
$ws.Range("U$row").Value = "NO"
$ws.Range("X$row").Value = "NO"
$ws.Range("AA$row").Value = 0

# Remaining row-29 columns (V, W, Y, Z, AB..AM) are blank in the source —
# give them the same empty text entry as the rest of the sheet.
$row29BlankCols = @("V", "W", "Y", "Z", "AB", "AC", "AD", "AE", "AF", "AG", "AH", "AI", "AJ", "AK", "AL", "AM")
foreach ($c in $row29BlankCols) {
    $cell = $ws.Range("$c$row")
    $cell.Value = "'"
    $cell.Style = "Normal"
}

